$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: turn into the new "smell overview" header row -----------------
# A1 keeps its bold header style, just gets new text.
$ws.Range("A1").Value = "Smell Name"

# B1 is brand-new, plain (no style).
$ws.Range("B1").Value = "Feature Envy"

# C1..L1 previously held smell-name headers with the "s=2" style; now they
# hold plain (unstyled) text.
$ws.Range("C1").Value = "Long method"
$ws.Range("D1").Value = "Inappropriate Intimacy"
$ws.Range("E1").Value = "Lazy Class"
$ws.Range("F1").Value = "Message Chain"
$ws.Range("G1").Value = "Many Parameters"
$ws.Range("H1").Value = "Duplicate Code"
$ws.Range("I1").Value = "Dead Code"
$ws.Range("J1").Value = "Unused Field"
$ws.Range("K1").Value = "No-op"
$ws.Range("L1").Value = "Use of Deprecated Interfaces"

$ws.Range("C1:L1").Style = "Normal"

# M1..AD1 (cols 13..30) lose their text but keep their "s=2" style (they
# become empty cells, matching the already-empty W1/AA1 gap cells).
for ($col = 13; $col -le 30; $col++) {
    $ws.Cells.Item(1, $col).ClearContents()
}

# --- Row 2: recreate the old header row content (used to be a fully empty
# row, now carries what row 1 used to hold: just "Program Name", bold). ----
$ws.Range("A2").Value = "Program Name"
$ws.Range("A2").Font.Bold = $true

# --- Row 3: add the feature-envy note next to "Afhaal Chinees". -----------
$ws.Range("B3").Value = "Yes, Sensor 2 is read in main program only to be fed into the my brick, would be better inside it"

# --- Column widths: new column B sized to fit the long note (target sheet
# stores width="75.42578125"; this host quantises ColumnWidth to 1/8-char
# steps, so 74.75 is the closest input that lands on the nearest
# representable width). -----------------------------------------------------
$ws.Range("B1").ColumnWidth = 74.75

# --- Selection, as recorded by the author when they saved the file. -------
$null = $ws.Range("B15").Select()
